$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet: ALC
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H112").Value = 84895.086
$ws.Range("J112").Value = 145089.58
$ws.Range("L112").Value = 435268.74
$ws.Range("N112").Value = -437484.74

$ws.Range("H116").Value = 15330
$ws.Range("I116").Value = 17995
$ws.Range("J116").Value = 10000
$ws.Range("K116").Value = 17995
$ws.Range("L116").Value = 10000
$ws.Range("M116").Value = -14553
$ws.Range("N116").Value = -16884

$ws.Range("H132").Value = 10102712
$ws.Range("I132").Value = 10754081
$ws.Range("J132").Value = 6500
$ws.Range("K132").Value = 32262243
$ws.Range("L132").Value = 19500
$ws.Range("M132").Value = -32259713
$ws.Range("N132").Value = -24560

$ws.Range("H141").Value = 1794.1818
$ws.Range("I141").Value = 1641.5238
$ws.Range("K141").Value = 4924.5714
$ws.Range("M141").Value = 255.4286000000002

# ---------------------------------------------------------------
# Sheet: ARM
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H32").Value = 4185.217
$ws.Range("I32").Value = 3192.1482
$ws.Range("J32").Value = 13122.833
$ws.Range("K32").Value = 3192.1482
$ws.Range("L32").Value = 13122.833
$ws.Range("M32").Value = -2905.1482
$ws.Range("N32").Value = -13696.833

$ws.Range("H45").Value = 6068.2
$ws.Range("I45").Value = 6037.1377
$ws.Range("K45").Value = 6037.1377
$ws.Range("M45").Value = -5660.1377

$ws.Range("H55").Value = 16999.334

$ws.Range("H61").Value = 3322
$ws.Range("I61").Value = 2631.2812
$ws.Range("J61").Value = 5777.8887
$ws.Range("K61").Value = 2631.2812
$ws.Range("L61").Value = 5777.8887
$ws.Range("M61").Value = -2419.2812
$ws.Range("N61").Value = -6201.8887

$ws.Range("H102").Value = 4475.5557
$ws.Range("I102").Value = 3843.1667
$ws.Range("K102").Value = 3843.1667
$ws.Range("M102").Value = -2221.1667

$ws.Range("H122").Value = 1895.8422
$ws.Range("I122").Value = 1622.1428
$ws.Range("J122").Value = 2662.2
$ws.Range("K122").Value = 4866.428400000001
$ws.Range("L122").Value = 7986.599999999999
$ws.Range("M122").Value = -2416.428400000001
$ws.Range("N122").Value = -12886.6

$ws.Range("H132").Value = 3810.5
$ws.Range("I132").Value = 3318.182
$ws.Range("J132").Value = 4893.6
$ws.Range("K132").Value = 9954.545999999998
$ws.Range("L132").Value = 14680.8
$ws.Range("M132").Value = -7424.545999999998
$ws.Range("N132").Value = -19740.8

$ws.Range("H134").Value = 89995
$ws.Range("J134").Value = 89995
$ws.Range("L134").Value = 89995
$ws.Range("N134").Value = -100135

$ws.Range("H135").Value = 96052.25
$ws.Range("J135").Value = 96052.25
$ws.Range("L135").Value = 96052.25
$ws.Range("N135").Value = -106192.25

$ws.Range("H136").Value = 3322
$ws.Range("I136").Value = 2631.2812
$ws.Range("J136").Value = 5777.8887
$ws.Range("K136").Value = 7893.8436
$ws.Range("L136").Value = 17333.6661
$ws.Range("M136").Value = -5343.8436
$ws.Range("N136").Value = -22433.6661

# ---------------------------------------------------------------
# Sheet: BSM
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H134").Value = 1584.2291
$ws.Range("I134").Value = 1648.4419
$ws.Range("K134").Value = 4945.3257
$ws.Range("M134").Value = -2410.3257

$ws.Range("H139").Value = 99967.664
$ws.Range("J139").Value = 99967.664
$ws.Range("L139").Value = 99967.664
$ws.Range("N139").Value = -110247.664

# ---------------------------------------------------------------
# Sheet: CRP
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H22").Value = 379.5
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

$ws.Range("H31").Value = 60273.832
$ws.Range("I31").Value = 68364
$ws.Range("J31").Value = 19823
$ws.Range("K31").Value = 68364
$ws.Range("L31").Value = 19823
$ws.Range("M31").Value = -68069
$ws.Range("N31").Value = -20413

$ws.Range("H34").Value = 60273.832
$ws.Range("I34").Value = 68364
$ws.Range("J34").Value = 19823
$ws.Range("K34").Value = 68364
$ws.Range("L34").Value = 19823
$ws.Range("M34").Value = -68162
$ws.Range("N34").Value = -20227

$ws.Range("H132").Value = 4364.5186
$ws.Range("I132").Value = 4289.68
$ws.Range("K132").Value = 12869.04
$ws.Range("M132").Value = -10339.04

# ---------------------------------------------------------------
# Sheet: CUL
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H5").Value = 50000
$ws.Range("J5").Value = 50000
$ws.Range("L5").Value = 150000
$ws.Range("N5").Value = -150224

$ws.Range("H113").Value = 668.0714
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 668.0714
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2004.2142
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -6344.2142

$ws.Range("H131").Value = 36831.723
$ws.Range("J131").Value = 2864.3333
$ws.Range("L131").Value = 8592.999899999999
$ws.Range("N131").Value = -18672.9999

$ws.Range("H135").Value = 50000
$ws.Range("J135").Value = 50000
$ws.Range("L135").Value = 450000
$ws.Range("N135").Value = -455070

$ws.Range("H138").Value = 35726420
$ws.Range("I138").Value = 71441850
$ws.Range("J138").Value = 10995.571
$ws.Range("K138").Value = 214325550
$ws.Range("L138").Value = 32986.713
$ws.Range("M138").Value = -214320410
$ws.Range("N138").Value = -43266.713

# ---------------------------------------------------------------
# Sheet: GSM
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H107").Value = 416.9375
$ws.Range("I107").Value = 342.63635
$ws.Range("J107").Value = 580.4
$ws.Range("K107").Value = 342.63635
$ws.Range("L107").Value = 580.4
$ws.Range("M107").Value = 1577.36365
$ws.Range("N107").Value = -4420.4

# ---------------------------------------------------------------
# Sheet: LTW
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H40").Value = 4812.222
$ws.Range("I40").Value = 3432.4614
$ws.Range("K40").Value = 3432.4614
$ws.Range("M40").Value = -3296.4614

$ws.Range("H100").Value = 2500
$ws.Range("I100").Value = 2000
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 2000
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -1459
$ws.Range("N100").Value = -4082

$ws.Range("H136").Value = 3969.4666
$ws.Range("I136").Value = 3817.926
$ws.Range("J136").Value = 5333.3335
$ws.Range("K136").Value = 11453.778
$ws.Range("L136").Value = 16000.0005
$ws.Range("M136").Value = -8903.778
$ws.Range("N136").Value = -21100.0005

$ws.Range("H138").Value = 146968
$ws.Range("J138").Value = 146968
$ws.Range("L138").Value = 146968
$ws.Range("N138").Value = -157248

$ws.Range("H140").Value = 164750
$ws.Range("I140").Value = 70000
$ws.Range("J140").Value = 449000
$ws.Range("K140").Value = 70000
$ws.Range("L140").Value = 449000
$ws.Range("M140").Value = -64820
$ws.Range("N140").Value = -459360

# ---------------------------------------------------------------
# Sheet: WVR
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H22").Value = 10370.333
$ws.Range("I22").Value = 1177.6666
$ws.Range("J22").Value = 14966.667
$ws.Range("K22").Value = 1177.6666
$ws.Range("L22").Value = 14966.667
$ws.Range("M22").Value = -884.6666
$ws.Range("N22").Value = -15552.667

$ws.Range("H23").Value = 1267.8334
$ws.Range("I23").Value = 1423.25
$ws.Range("J23").Value = 957
$ws.Range("K23").Value = 1423.25
$ws.Range("L23").Value = 957
$ws.Range("M23").Value = -1194.25
$ws.Range("N23").Value = -1415

$ws.Range("H39").Value = 9990.5
$ws.Range("I39").Value = 9990.5
$ws.Range("K39").Value = 9990.5
$ws.Range("M39").Value = -9577.5

$ws.Range("H136").Value = 2817.64
$ws.Range("J136").Value = 4377
$ws.Range("L136").Value = 13131
$ws.Range("N136").Value = -18231

Write-Output "All changes applied"
